# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (by copying the
#    "2021-Q4" sheet so its header row / column-A / header styling comes
#    along for free), then overwrite the fund rows with the new quarter's
#    holdings.
# 2. Prepend a new "2022-Q1" row to the "总计" (totals) summary sheet,
#    pushing the existing "2021-Q4" row down.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. Build the "2022-Q1" sheet from a copy of "2021-Q4" -----------------
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# Column B holds fund codes (e.g. "164811") and columns D:G hold plain
# decimal text (e.g. "0.23") in the source data - both look numeric, so
# Excel would silently coerce them to numbers on assignment. Force the
# range to Text first so the values stick as strings, then clear the
# number-format back off again (this keeps the cell style itself at the
# workbook default - just like the source sheet - while the stored value
# remains text).
$dataText = $q1.Range("B2:G4")
$dataText.NumberFormat = "@"

$q1.Range("B2").Value = "164811"
$q1.Range("C2").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
$q1.Range("D2").Value = "0.23"
$q1.Range("E2").Value = "94.28"
$q1.Range("F2").Value = "4.93"
$q1.Range("G2").Value = "0.0113"

$q1.Range("B3").Value = "512780"
$q1.Range("C3").Value = "广发中证京津冀协同发展主题ETF"
$q1.Range("D3").Value = "0.13"
$q1.Range("E3").Value = "98.52"
$q1.Range("F3").Value = "3.46"
$q1.Range("G3").Value = "0.0045"

$q1.Range("B4").Value = "164825"
$q1.Range("C4").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$q1.Range("D4").Value = "0.06"
$q1.Range("E4").Value = "94.28"
$q1.Range("F4").Value = "4.93"
$q1.Range("G4").Value = "0.0030"

$dataText.ClearFormats()

$q1.Range("H2").Value = 1
$q1.Range("H3").Value = 2
$q1.Range("H4").Value = 1

# --- 2. Prepend the "2022-Q1" row on the "总计" sheet ----------------------
$total = $wb.Worksheets.Item("总计")

$oldDate = $total.Range("B2").Value()
$oldCount = $total.Range("C2").Value()
$oldMarket = $total.Range("D2").Value()

# Copy row 2's formatting down to row 3 before overwriting row 2, so the
# bold/bordered style on column A follows the shifted "2021-Q4" row.
$total.Range("A2:D2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldDate
$total.Range("C3").Value = $oldCount
$total.Range("D3").Value = $oldMarket

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.02
